$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset styles to Normal FIRST for all cells we will touch, so any later
# quote-prefixed text values reuse the same minimal style delta.
$ws.Range("A2:G5").Style = "Normal"

$ws.Range("A2").Value = 495
$ws.Range("B2").Value = "simon"
$ws.Range("C2").Value = "'9"
$ws.Range("D2").Value = "perro"
$ws.Range("E2").Value = "es mi perrito"
$ws.Range("F2").Value = "esta muy bien"
$ws.Range("G2").Value = $true

$ws.Range("A3").Value = 546
$ws.Range("B3").Value = "max"
$ws.Range("C3").Value = "'2"
$ws.Range("D3").Value = "es perro"
$ws.Range("E3").Value = "lo quiero mucho"
$ws.Range("F3").Value = "le falta un ojo"
$ws.Range("G3").Value = $false

$ws.Range("A4").Value = 88
$ws.Range("B4").Value = "pepe"
$ws.Range("C4").Value = "'3"
$ws.Range("D4").Value = "reptil"
$ws.Range("E4").Value = "es una tortuga"
$ws.Range("F4").Value = "esta bien"
$ws.Range("G4").Value = $false

$ws.Range("A5").Value = 61
$ws.Range("B5").Value = "'"
$ws.Range("C5").Value = "'4"
$ws.Range("D5").Value = "'"
$ws.Range("E5").Value = "'"
$ws.Range("F5").Value = "'"
$ws.Range("G5").Value = $false

# Re-apply Normal to strip the transient quote-prefix style now that every
# value has been committed.
$ws.Range("A2:G5").Style = "Normal"

$ws.Rows.Item(6).Delete()

$pct = $wb.Styles.Item("Porcentaje")
$pct.Delete()

[void]$ws.Range("K7").Select()
Write-Host "edits applied"
